$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(351.0, "Lucas", "Agricultor"),
    @(521.0, "Ana", "Engenheira"),
    @(572.0, "Ivone", "Violeira"),
    @(701.0, "Bruna", "Sanfoneira")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $obj = $data[$i]
    $ws.Cells.Item($row, 1).Value = [double]$obj[0]
    $ws.Cells.Item($row, 2).Value = $obj[1]
    $ws.Cells.Item($row, 3).Value = $obj[2]
}
